$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 334, shifting existing rows 334:406 down to 335:407
$ws.Rows.Item(334).Insert()

# New row 334 mirrors the row layout (Mercado ID ... Clasificacion) with new data
$ws.Cells.Item(334, 1).Value = 3
$ws.Cells.Item(334, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(334, 3).Value = "Coquimbo"
$ws.Cells.Item(334, 4).Value = [DateTime]"2022-08-25"
$ws.Cells.Item(334, 5).Value = 5
$ws.Cells.Item(334, 6).Value = 100112043
$ws.Cells.Item(334, 7).Value = "Pepino ensalada"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 100
$ws.Cells.Item(334, 11).Value = 21000
$ws.Cells.Item(334, 12).Value = 22000
$ws.Cells.Item(334, 13).Value = 21550
$ws.Cells.Item(334, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(334, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(334, 16).Value = 308
$ws.Cells.Item(334, 17).Value = 70
$ws.Cells.Item(334, 18).Value = "Hortaliza"
